$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.751760721206665
$ws.Range("B1").Value = 2.644430875778198
$ws.Range("C1").Value = 3.271781206130981
$ws.Range("D1").Value = 1.248926401138306
$ws.Range("E1").Value = 0.8306245803833008
